# Update cryptos list prices and volume percentages (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").Value = "30.346.57"
    $ws.Range("E2").Value = "  +0.43%  "
    $ws.Range("D3").Value = "1.937.39"
    $ws.Range("E3").Value = "  +0.65%  "
    $ws.Range("E4").Value = "  -0.26%  "
    $ws.Range("D5").Value = "252.02"
    $ws.Range("E5").Value = "  +2.00%  "
    $ws.Range("D6").Value = "0.7240"
    $ws.Range("E6").Value = "  +3.68%  "
    $ws.Range("E7").Value = "  -0.24%  "
    $ws.Range("D8").Value = "0.3318"
    $ws.Range("E8").Value = "  +2.26%  "
    $ws.Range("D9").Value = "28.00"
    $ws.Range("E9").Value = "  +6.05%  "
    $ws.Range("D10").Value = "0.07271"
    $ws.Range("E10").Value = "  +6.60%  "
    $ws.Range("E11").Value = "  +1.84%  "
    $ws.Range("D12").Value = "0.08105"
    $ws.Range("E12").Value = "  +1.99%  "
    $ws.Range("D13").Value = "1.935.74"
    $ws.Range("E13").Value = "  +0.52%  "
    $ws.Range("D14").Value = "5.490"
    $ws.Range("E14").Value = "  +2.00%  "
    $ws.Range("D15").Value = "95.01"
    $ws.Range("E15").Value = "  +1.11%  "
    $ws.Range("D16").Value = "15.16"
    $ws.Range("E16").Value = "  +4.72%  "
    $ws.Range("D17").Value = "30.343.38"
    $ws.Range("E17").Value = "  +0.34%  "
    $ws.Range("D18").Value = "0.000008253"
    $ws.Range("E18").Value = "  +5.48%  "
    $ws.Range("D19").Value = "253.50"
    $ws.Range("E19").Value = "  -2.36%  "
    $ws.Range("D20").Value = "5.843"
    $ws.Range("E20").Value = "  +0.32%  "
    $ws.Range("E21").Value = "  +0.58%  "
    $ws.Range("E22").Value = "  -0.23%  "
    $ws.Range("D23").Value = "1.000"
    $ws.Range("E23").Value = "  -0.20%  "
    $ws.Range("D24").Value = "6.967"
    $ws.Range("E24").Value = "  +2.37%  "
    $ws.Range("D25").Value = "9.778"
    $ws.Range("E25").Value = "  +1.73%  "
    $ws.Range("D26").Value = "165.99"
    $ws.Range("E26").Value = "  +4.34%  "
    $ws.Range("D27").Value = "2.350"
    $ws.Range("E27").Value = "  +5.75%  "
    $ws.Range("D28").Value = "19.36"
    $ws.Range("E28").Value = "  +3.19%  "
    $ws.Range("D29").Value = "0.1302"
    $ws.Range("E29").Value = "  -0.87%  "
    $ws.Range("E30").Value = "  +0.93%  "
    $ws.Range("D31").Value = "1.543"
    $ws.Range("E31").Value = "  -0.71%  "
    $ws.Range("D32").Value = "4.445"
    $ws.Range("E32").Value = "  +1.12%  "
    $ws.Range("D33").Value = "4.218"
    $ws.Range("E33").Value = "  +0.84%  "
    $ws.Range("D34").Value = "0.05259"
    $ws.Range("E34").Value = "  +4.50%  "
    $ws.Range("D35").Value = "1.271"
    $ws.Range("E35").Value = "  +6.69%  "
    $ws.Range("D36").Value = "0.7514"
    $ws.Range("E36").Value = "  +0.74%  "
    $ws.Range("D37").Value = "2.770"
    $ws.Range("E37").Value = "  +2.27%  "
    $ws.Range("D38").Value = "0.01977"
    $ws.Range("E38").Value = "  +3.10%  "
    $ws.Range("D39").Value = "2.804"
    $ws.Range("E39").Value = "  +0.23%  "
    $ws.Range("D40").Value = "79.61"
    $ws.Range("E40").Value = "  -0.43%  "
    $ws.Range("D41").Value = "6.453"
    $ws.Range("E41").Value = "  -0.84%  "
    $ws.Range("D42").Value = "0.4562"
    $ws.Range("E42").Value = "  +3.76%  "
    $ws.Range("D43").Value = "2.035"
    $ws.Range("E43").Value = "  +0.34%  "
    $ws.Range("D44").Value = "0.8441"
    $ws.Range("E44").Value = "  +1.29%  "
    $ws.Range("D45").Value = "1.000"
    $ws.Range("E45").Value = "  -0.18%  "
    $ws.Range("D46").Value = "102.08"
    $ws.Range("E46").Value = "  +0.36%  "
    $ws.Range("D47").Value = "9.795"
    $ws.Range("E47").Value = "  +1.72%  "
    $ws.Range("D48").Value = "7.462"
    $ws.Range("E48").Value = "  +3.75%  "
    $ws.Range("D50").Value = "0.4210"
    $ws.Range("E50").Value = "  +3.84%  "
    $ws.Range("D51").Value = "0.06048"
    $ws.Range("E51").Value = "  +1.78%  "

Write-Output "Updated cryptos list"
